$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16
$ws.Range("D16").Value = "image_20250807111026_ppp0.jpg"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "641,529,688,576"

# Row 17
$ws.Range("D17").Value = "image_20250807111026_ppp0.jpg"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "793,481,831,526"

# Row 18
$ws.Range("D18").Value = "image_20250808221835_ppp0.jpg"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "1182,405,1231,455"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "0.76"
